$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that hold numeric-looking text (e.g. "1.000", "0.5164")
# must be forced to Text format first, otherwise Excel auto-converts the
# string into a number when the value is assigned.

$ws.Range("D2").Value = "30.234.36"
$ws.Range("E2").Value = "  +3.35%  "

$ws.Range("D3").Value = "1.897.10"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("E4").Value = "  -0.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.37"
$ws.Range("E5").Value = "  +3.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5164"
$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4012"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08451"
$ws.Range("E9").Value = "  +0.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.117"
$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.30"
$ws.Range("E12").Value = "  +12.86%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.429"
$ws.Range("E13").Value = "  +2.73%  "

$ws.Range("D14").Value = "1.892.78"
$ws.Range("E14").Value = "  +0.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.342"
$ws.Range("E15").Value = "  +0.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.78"
$ws.Range("E17").Value = "  +1.80%  "

$ws.Range("E18").Value = "  +0.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06650"
$ws.Range("E19").Value = "  -1.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.25"
$ws.Range("E20").Value = "  +2.39%  "

$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.951"
$ws.Range("E22").Value = "  -0.97%  "

$ws.Range("D23").Value = "30.237.71"
$ws.Range("E23").Value = "  +3.36%  "

$ws.Range("E24").Value = "  +1.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.227"
$ws.Range("E25").Value = "  +0.61%  "

$ws.Range("D26").Value = "2.109.25"
$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.67"
$ws.Range("E27").Value = "  +3.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.09"
$ws.Range("E28").Value = "  +1.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.365"
$ws.Range("E29").Value = "  -2.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.93"
$ws.Range("E30").Value = "  +0.98%  "

$ws.Range("E31").Value = "  +3.68%  "

$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.062"
$ws.Range("E33").Value = "  -1.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.758"
$ws.Range("E34").Value = "  +2.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02493"
$ws.Range("E35").Value = "  +0.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06559"
$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.258"
$ws.Range("E37").Value = "  +1.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2203"
$ws.Range("E38").Value = "  +0.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.219"
$ws.Range("E39").Value = "  -0.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.81"
$ws.Range("E40").Value = "  +4.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6503"
$ws.Range("E41").Value = "  +0.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.716"
$ws.Range("E42").Value = "  -3.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.234"
$ws.Range("E43").Value = "  +0.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6109"
$ws.Range("E44").Value = "  +1.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.21"
$ws.Range("E45").Value = "  +0.56%  "

$ws.Range("E46").Value = "  +0.80%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.058"
$ws.Range("E47").Value = "  +0.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.235"
$ws.Range("E48").Value = "  +0.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.40"
$ws.Range("E49").Value = "  +0.83%  "

$ws.Range("E50").Value = "  -0.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.96"
$ws.Range("E51").Value = "  +2.03%  "
